$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Last status check on" timestamp in F1 (09:00 -> 09:15)
$ws.Range("F1").Value = "Last status check on: 15.02.2022 09:15"

# 2. Row 8 (Benzina Albert Modrice): D8 and E8 were written as text by the
#    scraper; convert them to real numeric/date values like the other rows.
$ws.Range("D8").Value = -0.7

# E8 = 2022-02-15 09:03:07 expressed as an Excel serial date/time, formatted
# the same way as the other "Old Datum" cells (style index 2 -> numFmtId 165).
$ws.Range("E8").Value = 44607.37716435185
$ws.Range("E8").NumberFormat = $ws.Range("E2").NumberFormat
